$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.816.53"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.438.35"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +5.86%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "484.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +14.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +8.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.461.89"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0970"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +11.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.84%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.326"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +8.71%  "

# Row 13
$ws.Range("E13").Value = "  +1.44%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.857.43"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.938.83"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.79%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.59"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +9.98%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +14.61%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.452.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +11.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.91"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +10.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.995"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.32%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.30%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +13.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.406"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +11.94%  "

# Row 27
$ws.Range("E27").Value = "  +0.68%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.574.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.07%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0781"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +18.57%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.33%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.01"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.33%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.76%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.14"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +9.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.27%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.54%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.857"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.83"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.76%  "

# Row 40
$ws.Range("E40").Value = "  +0.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.601"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.31%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0547"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.16%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +10.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.30"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +11.26%  "

# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.68"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +11.80%  "

# Row 46
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.08"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.63%  "

# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "254.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +28.34%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0901"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +10.11%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.948.42"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.03%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0223"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.24"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +11.13%  "
